# 301312-智立方: add 2022-Q4 data.
#
# Before: two sheets, "总计" (running totals) and "2022-Q3" (that quarter's
# fund holdings).
# After:  three sheets - "总计", "2022-Q4" (new fund holdings), and
# "2022-Q3" (the original fund-holdings data, preserved verbatim under its
# original tab).

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$qSheet = $wb.Worksheets.Item("2022-Q3")

# Duplicate the current "2022-Q3" sheet, placing the copy right after it, so
# the old quarter's numbers stay intact on their own (still so-named) tab.
$qSheet.Copy($null, $qSheet)
$copySheet = $wb.Worksheets.Item(3)

# Free up the "2022-Q3" name on the original sheet before claiming it again
# for the copy.
$qSheet.Name = "2022-Q4"
$copySheet.Name = "2022-Q3"

# Fill in the 2022-Q4 fund figures on the renamed sheet.
$qSheet.Range("D2").Value = "'0.42"
$qSheet.Range("E2").Value = "'92.60"
$qSheet.Range("F2").Value = "'4.46"
$qSheet.Range("G2").Value = "'0.0187"
$qSheet.Range("H2").Value = 6

# The leading apostrophes above force D2:G2 to be stored as text (matching
# the source data, which keeps these figures as text), but they also stamp
# a "quote prefix" onto the cell format. Clear that incidental formatting so
# the cells stay plain/unstyled, same as the rest of the data row.
$qSheet.Range("D2:G2").ClearFormats()

# The header row and the leading A2 marker cell pick up the same style used
# for the header/marker cells on the "总计" sheet.
$totalSheet.Range("B1").Copy()
$qSheet.Range("B1:H1").PasteSpecial(-4122)
$totalSheet.Range("A2").Copy()
$qSheet.Range("A2").PasteSpecial(-4122)

# Update "总计": the existing row now reports 2022-Q4, and a new row is
# appended for 2022-Q3 (same counts/value, same row styling).
$totalSheet.Range("B2").Value = "2022-Q4"

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.02
